$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "Hello Canada"
$ws.Range("A2").Select()
